# Update NATMI LR-pair results (Agrn-Musk) with newly recomputed TPM-based
# statistics. The new run only reports target clusters FAPs/MuSCs (rows
# with Target cluster = ECs are no longer produced), shrinking the table
# from 9 data rows (rows 2-10) to 6 data rows (rows 2-7), and every
# remaining numeric column (G:T) is recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the three trailing rows (old "Sending cluster = MuSCs" x all targets
# block is gone; overall the table is now 6 data rows instead of 9).
$ws.Range("A8:T10").EntireRow.Delete()

# Row 2: ECs -> FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Agrn"
$ws.Range("C2").Value = "Musk"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 8.382531
$ws.Range("H2").Value = 25.147593
$ws.Range("I2").Value = 0.3278601051951505
$ws.Range("J2").Value = 0.3278601051951506
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 8.908863
$ws.Range("N2").Value = 26.726589
$ws.Range("O2").Value = 0.6232066589444157
$ws.Range("P2").Value = 0.6232066589444157
$ws.Range("Q2").Value = 74.67882027225301
$ws.Range("R2").Value = 672.1093824502771
$ws.Range("S2").Value = 0.2043246007598344
$ws.Range("T2").Value = 0.2043246007598345

# Row 3: ECs -> MuSCs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Agrn"
$ws.Range("C3").Value = "Musk"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 8.382531
$ws.Range("H3").Value = 25.147593
$ws.Range("I3").Value = 0.3278601051951505
$ws.Range("J3").Value = 0.3278601051951506
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.386335666666667
$ws.Range("N3").Value = 16.159007
$ws.Range("O3").Value = 0.3767933410555842
$ws.Range("P3").Value = 0.3767933410555843
$ws.Range("Q3").Value = 45.151125702239
$ws.Range("R3").Value = 406.360131320151
$ws.Range("S3").Value = 0.1235355044353161
$ws.Range("T3").Value = 0.1235355044353161

# Row 4: FAPs -> FAPs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Agrn"
$ws.Range("C4").Value = "Musk"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.399531333333333
$ws.Range("H4").Value = 19.198594
$ws.Range("I4").Value = 0.2503004183517279
$ws.Range("J4").Value = 0.250300418351728
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 8.908863
$ws.Range("N4").Value = 26.726589
$ws.Range("O4").Value = 0.6232066589444157
$ws.Range("P4").Value = 0.6232066589444157
$ws.Range("Q4").Value = 57.012547912874
$ws.Range("R4").Value = 513.112931215866
$ws.Range("S4").Value = 0.1559888874533699
$ws.Range("T4").Value = 0.1559888874533699

# Row 5: FAPs -> MuSCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Agrn"
$ws.Range("C5").Value = "Musk"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.399531333333333
$ws.Range("H5").Value = 19.198594
$ws.Range("I5").Value = 0.2503004183517279
$ws.Range("J5").Value = 0.250300418351728
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.386335666666667
$ws.Range("N5").Value = 16.159007
$ws.Range("O5").Value = 0.3767933410555842
$ws.Range("P5").Value = 0.3767933410555843
$ws.Range("Q5").Value = 34.47002387068422
$ws.Range("R5").Value = 310.230214836158
$ws.Range("S5").Value = 0.09431153089835803
$ws.Range("T5").Value = 0.09431153089835807

# Row 6: MuSCs -> FAPs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Agrn"
$ws.Range("C6").Value = "Musk"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 10.78533933333333
$ws.Range("H6").Value = 32.356018
$ws.Range("I6").Value = 0.4218394764531215
$ws.Range("J6").Value = 0.4218394764531215
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 8.908863
$ws.Range("N6").Value = 26.726589
$ws.Range("O6").Value = 0.6232066589444157
$ws.Range("P6").Value = 0.6232066589444157
$ws.Range("Q6").Value = 96.085110529178
$ws.Range("R6").Value = 864.765994762602
$ws.Range("S6").Value = 0.2628931707312114
$ws.Range("T6").Value = 0.2628931707312114

# Row 7: MuSCs -> MuSCs
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Agrn"
$ws.Range("C7").Value = "Musk"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 10.78533933333333
$ws.Range("H7").Value = 32.356018
$ws.Range("I7").Value = 0.4218394764531215
$ws.Range("J7").Value = 0.4218394764531215
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.386335666666667
$ws.Range("N7").Value = 16.159007
$ws.Range("O7").Value = 0.3767933410555842
$ws.Range("P7").Value = 0.3767933410555843
$ws.Range("Q7").Value = 58.09345792823622
$ws.Range("R7").Value = 522.8411213541259
$ws.Range("S7").Value = 0.1589463057219101
$ws.Range("T7").Value = 0.1589463057219101
